# Updated cryptos list values (price + 1h volume change) per target diff.
# For Price (column D) cells whose new text would otherwise be auto-parsed
# as a plain number by Excel's smart cell-entry, force the cell to Text
# format first so the written value stays a string (matching the source
# workbook, which stores every Price/Volume cell as text).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.205.09"
$ws.Range("E2").Value = "  -0.66%  "

$ws.Range("D3").Value = "1.829.90"
$ws.Range("E3").Value = "  -0.73%  "

$ws.Range("E4").Value = "  +0.23%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "237.39"
$ws.Range("E5").Value = "  -1.17%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6085"
$ws.Range("E6").Value = "  -3.78%  "

$ws.Range("E7").Value = "  +0.22%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07098"
$ws.Range("E8").Value = "  -5.03%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2818"
$ws.Range("E9").Value = "  -3.01%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "23.86"
$ws.Range("E10").Value = "  -4.63%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07643"
$ws.Range("E11").Value = "  -1.28%  "

$ws.Range("D12").Value = "1.843.87"
$ws.Range("E12").Value = "  -0.13%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.811"
$ws.Range("E13").Value = "  -3.56%  "

$ws.Range("D16").Value = "2.067.11"
$ws.Range("E16").Value = "  -1.34%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "79.57"
$ws.Range("E17").Value = "  -2.97%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "5.964"
$ws.Range("E18").Value = "  -4.88%  "

$ws.Range("D19").Value = "29.208.19"

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "229.21"
$ws.Range("E20").Value = "  -0.30%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.035"
$ws.Range("E23").Value = "  -5.17%  "

$ws.Range("E24").Value = "  +0.21%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "155.56"
$ws.Range("E25").Value = "  -1.75%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.108"
$ws.Range("E26").Value = "  -4.57%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1302"
$ws.Range("E27").Value = "  -4.00%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "16.73"
$ws.Range("E28").Value = "  -4.22%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.06783"
$ws.Range("E29").Value = "  +3.57%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.477"
$ws.Range("E30").Value = "  +3.45%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.457"
$ws.Range("E31").Value = "  -1.92%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.837"
$ws.Range("E32").Value = "  -5.88%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.840"
$ws.Range("E33").Value = "  -5.22%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.131"
$ws.Range("E34").Value = "  -0.92%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.733"
$ws.Range("E35").Value = "  -5.79%  "

$ws.Range("E36").Value = "  -6.20%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.557"
$ws.Range("E37").Value = "  -0.76%  "

$ws.Range("D38").Value = "1.231.02"
$ws.Range("E38").Value = "  -1.58%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.764"
$ws.Range("E39").Value = "  -1.87%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01769"
$ws.Range("E40").Value = "  -4.81%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.602"
$ws.Range("E41").Value = "  -2.55%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9209"
$ws.Range("E42").Value = "  -1.48%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.002"
$ws.Range("E43").Value = "  +0.17%  "

$ws.Range("D44").Value = "1.984.28"
$ws.Range("E44").Value = "  -1.93%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "100.94"
$ws.Range("E45").Value = "  -0.02%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "63.49"
$ws.Range("E46").Value = "  -3.04%  "

$ws.Range("E47").Value = "  -1.37%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.631"
$ws.Range("E48").Value = "  -5.06%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.578"
$ws.Range("E49").Value = "  -4.82%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.1086"
$ws.Range("E50").Value = "  -5.46%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.533"
$ws.Range("E51").Value = "  -7.49%  "

$ws.Range("B14").Value = "ShibaInu"
$ws.Range("C14").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.00001001"
$ws.Range("E14").Value = "  -2.38%  "

$ws.Range("B15").Value = "Polygon"
$ws.Range("C15").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6329"
$ws.Range("E15").Value = "  -6.74%  "

$ws.Range("B21").Value = "Avalanche"
$ws.Range("C21").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.81"
$ws.Range("E21").Value = "  -4.19%  "

$ws.Range("B22").Value = "Dai"
$ws.Range("C22").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.002"
$ws.Range("E22").Value = "  +0.25%  "
